# Airbnb task-tracker workbook: "fix bug and done XacNhan"
#
# - Row 11 ("Chi tiet phong o + Danh sach tien ich"): assign the owner
#   name ("Anh Ngoc") in the Names column.
# - Row 12 ("Xac nhan va dat cho + Dang xuat"): mark the build as finished
#   -> Finish Build date + % Build + owner name (Anh Ngoc).
# - Move the active selection to G11 (where the editor was last working).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: fill in the "Names" cell (E11) with the builder's name ----
# Re-use the formatting already used by the other filled-in Name cells
# in this block (E8/E9/E10) so the new cell matches the table styling.
$ws.Range("E9").Copy()
$ws.Range("E11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E11").Value = "Ánh Ngọc"

# --- Row 12: Finish-Build date, %-Build, and Names ----------------------
# Finish Build (C12) takes the same formatting as the cell above it (C11).
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C12").Value = 44642

# % Build (D12) keeps its existing number format; just record 60% done.
$ws.Range("D12").Value = 0.6

# Names (E12) same treatment as E11 above.
$ws.Range("E9").Copy()
$ws.Range("E12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E12").Value = "Ánh Ngọc"

$excel.CutCopyMode = 0

# --- Restore a clean selection, matching where the author left off -----
$ws.Range("G11").Select() | Out-Null
